# repull data, push all data, mean calculation
# Update the dSF (column F) values (and, for row 24, dS0/IP/I0 too) to reflect
# the repulled/recomputed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell -> new value
$updates = @{
    "F2"  = -6
    "F3"  = 3
    "F5"  = 12
    "F6"  = -3
    "F7"  = -9
    "F10" = -5
    "F12" = 1
    "F13" = -10
    "F14" = -5
    "F19" = -3
    "F21" = -2
    "F22" = -7
    "E24" = -2
    "F24" = -1
    "H24" = 2
    "I24" = 6
    "F25" = 5
    "F26" = -1
    "F30" = -2
    "F31" = 1
    "F32" = -6
    "F35" = -5
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
